$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Unneeded review points": wipe out the two review entries that had been
# filled in (CR_REV_0001 / CR_REV_0002 and their Reviewer/Version/Section/
# Problem/Solution columns), leaving just the Status ("Open") like the
# still-blank rows below them.
$ws.Range("A2:G2").ClearContents()
$ws.Range("A3:G3").ClearContents()

# Row 3 (and row 4, which was already blank) revert from the "filled-in"
# look back to the plain style used by the rest of the blank rows.
$ws.Range("A3").Font.Name = "Arial"
$ws.Range("C3:G3").Font.Name = "Arial"
$ws.Range("A4").Font.Name = "Arial"
$ws.Range("C4:G4").Font.Name = "Arial"

# "revert changes": drop the two extra blank review rows that had been
# appended, going back to the original, shorter table.
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(8).Delete()

# Shrink the conditional-formatting rules' ranges back down to match the
# smaller table instead of leaving them sized for the old 9-row layout.
$full = $ws.Range("A2:H9")
$fcs = $full.FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    $fc.ModifyAppliesToRange($ws.Range("A2:H7"))
}

# Restore the previously-selected row.
[void]$ws.Rows.Item(2).Select()
